$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(308051846, Eyal  Sofer: 4,9)"
$ws.Range("B1").Value = "(312049950, Molham  Peretz: -1,-3)"
$ws.Range("C1").Value = "(308073899, Anan  Kirshenbaum: -2,-9)"
$ws.Range("D1").Value = "(318869187, Soaad  Leibovich: -9,-10)"
$ws.Range("E1").Value = "(205898513, Asaf  Braymok: -3,-7)"
$ws.Range("F1").Value = "(318428158, Tal  Asulin: -1,-8)"
$ws.Range("G1").Value = "(316028364, Sami  Castro: 2,-10)"

$ws.Range("A3").Value = "cost: 218.5182676381342"
$ws.Range("A4").Value = "time: 38.703653527626834"
